$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Values in column D that are plain numeric strings (e.g. "214.62") would be
# auto-coerced to numbers by Excel on assignment; the source data keeps them as
# literal text (e.g. "1.001" meaning $1.001, not the number 1.001, and trailing
# zeros such as "0.5076" must survive). Force those specific cells to Text format
# first so the assigned string is stored verbatim.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

$ws.Range("D2").Value = "25.915.36"
$ws.Range("E2").Value = "  -0.15%  "
$ws.Range("D3").Value = "1.638.39"
$ws.Range("E4").Value = "  -0.66%  "
$ws.Range("D5").Value = "214.62"
$ws.Range("E5").Value = "  -0.40%  "
$ws.Range("D6").Value = "0.5076"
$ws.Range("E6").Value = "  -0.12%  "
$ws.Range("E7").Value = "  -0.21%  "
$ws.Range("E8").Value = "  +0.25%  "
$ws.Range("D9").Value = "0.06366"
$ws.Range("E9").Value = "  -0.31%  "
$ws.Range("E10").Value = "  +1.39%  "
$ws.Range("D11").Value = "0.07747"
$ws.Range("E11").Value = "  -0.53%  "
$ws.Range("D12").Value = "4.298"
$ws.Range("E12").Value = "  -0.20%  "
$ws.Range("D13").Value = "1.632.74"
$ws.Range("E13").Value = "  -0.54%  "
$ws.Range("D14").Value = "0.5466"
$ws.Range("E14").Value = "  -0.04%  "
$ws.Range("E15").Value = "  -1.46%  "
$ws.Range("D16").Value = "64.21"
$ws.Range("E16").Value = "  -0.39%  "
$ws.Range("D17").Value = "25.933.96"
$ws.Range("E17").Value = "  -0.23%  "
$ws.Range("D18").Value = "1.001"
$ws.Range("E18").Value = "  -0.27%  "
$ws.Range("E19").Value = "  +0.62%  "
$ws.Range("D20").Value = "195.97"
$ws.Range("E20").Value = "  -1.11%  "
$ws.Range("D21").Value = "9.946"
$ws.Range("E21").Value = "  -0.18%  "
$ws.Range("D22").Value = "6.137"
$ws.Range("E22").Value = "  +1.43%  "
$ws.Range("E23").Value = "  -0.31%  "
$ws.Range("D24").Value = "1.895"
$ws.Range("E24").Value = "  +0.68%  "
$ws.Range("D25").Value = "143.33"
$ws.Range("E25").Value = "  +1.24%  "
$ws.Range("D26").Value = "0.1259"
$ws.Range("D27").Value = "6.845"
$ws.Range("E27").Value = "  -0.58%  "
$ws.Range("D28").Value = "15.64"
$ws.Range("E28").Value = "  -0.69%  "
$ws.Range("D29").Value = "1.236"
$ws.Range("E29").Value = "  -0.07%  "
$ws.Range("D30").Value = "0.04886"
$ws.Range("E30").Value = "  -3.19%  "
$ws.Range("D31").Value = "3.246"
$ws.Range("E31").Value = "  -0.47%  "
$ws.Range("D32").Value = "3.206"
$ws.Range("E32").Value = "  +0.38%  "
$ws.Range("D33").Value = "1.554"
$ws.Range("E33").Value = "  +0.80%  "
$ws.Range("D34").Value = "2.374"
$ws.Range("E34").Value = "  +0.49%  "
$ws.Range("D35").Value = "0.9141"
$ws.Range("E35").Value = "  +2.06%  "
$ws.Range("D36").Value = "2.571"
$ws.Range("E36").Value = "  -0.98%  "
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").Value = "0.5527"
$ws.Range("E37").Value = "  +0.63%  "
$ws.Range("B38").Value = "Maker"
$ws.Range("C38").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D38").Value = "1.126.06"
$ws.Range("E38").Value = "  -0.60%  "
$ws.Range("E39").Value = "  +0.42%  "
$ws.Range("D40").Value = "1.002"
$ws.Range("E40").Value = "  -0.27%  "
$ws.Range("D41").Value = "5.603"
$ws.Range("E41").Value = "  -0.43%  "
$ws.Range("D42").Value = "0.8041"
$ws.Range("E42").Value = "  -1.64%  "
$ws.Range("D43").Value = "98.53"
$ws.Range("D44").Value = "0.0₈120"
$ws.Range("E44").Value = "  -9.89%  "
$ws.Range("D45").Value = "1.772.50"
$ws.Range("E45").Value = "  -0.34%  "
$ws.Range("D46").Value = "0.4487"
$ws.Range("E46").Value = "  -0.93%  "
$ws.Range("D47").Value = "55.19"
$ws.Range("E47").Value = "  +0.48%  "
$ws.Range("E48").Value = "  -0.09%  "
$ws.Range("D49").Value = "0.05185"
$ws.Range("E49").Value = "  +2.18%  "
$ws.Range("D50").Value = "7.498"
$ws.Range("E50").Value = "  +0.99%  "
$ws.Range("E51").Value = "  -0.16%  "
